$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the old_vars/new_vars mapping table ---
# Cells are written in the same order the new label text first appears in
# the final sheet (D1, A5, B5, A4, A6) so the regenerated shared-strings
# table lands in the same order as the target workbook.

# Row 1: header row; D1 label changes (expectedwear -> sensor_lifetime)
$ws.Range("D1").Value = "sensor_lifetime"

# Row 5: Historic Glucose mmol/L/sensorglucose -> Record Type/scan_yn (new row)
$ws.Range("A5").Value = "Record Type"
$ws.Range("B5").Value = "scan_yn"

# Row 4: Device Timestamp/timestamp -> Historic Glucose(mmol/L)/sensorglucose
$ws.Range("A4").Value = "Historic Glucose(mmol/L)"
$ws.Range("B4").Value = "sensorglucose"

# Row 6: Scan Glucose mmol/L -> Scan Glucose(mmol/L) (space before parens removed)
$ws.Range("A6").Value = "Scan Glucose(mmol/L)"
$ws.Range("B6").Value = "scanglucose"

# Row 2: Device/device -> Serial Number/id  (the old "Device" row is dropped,
# rows shift up by one)
$ws.Range("A2").Value = "Serial Number"
$ws.Range("B2").Value = "id"

# Row 3: Serial Number/id -> Device Timestamp/timestamp
$ws.Range("A3").Value = "Device Timestamp"
$ws.Range("B3").Value = "timestamp"

# Rows 7-13 (dexcomg6 / other groups) keep the same text, only the
# underlying shared-string indices shift because of the table edits above -
# no cell-value change required there.

# --- View/selection metadata ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E9").Select()
$excel.ActiveWindow.Zoom = 100
